# Apply the "added columns, flags and comp description" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data table (A1:F8)
$data = @(
    @("Team Name", "Country", "Running order", "Day score", "Time ranking", "Total score"),
    @("IceBerg", "Canada", 1, 50, 1, 100),
    @("Team Aritra", "India", 4, 95, 4, 130),
    @("Navier USN", "Norway", 3, 99, 7, 120),
    @("Vortex", "Norway", 5, 98, 3, 170),
    @("AGH Solar Boat", "Poland", 7, 95, 5, 200),
    @("Técnico Solar Boat", "Portugal", 6, 90, 2, 140),
    @("StrathVoyager", "Scotland", 2, 60, 6, 95)
)

# Clear out the old C column header ("Score") / stray formatting from the
# previous narrower table before writing the new, wider one.
$ws.Range("A1:F8").Clear()

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Column widths (bestFit, customWidth) for Team Name (A) and Running order (C)
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(3).ColumnWidth = 10.8

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moved to G9
$ws.Range("G9").Select()
